# Auto-generated edit script: updates profit/price calculation cells
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled runner refresh of market price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 742.8461
$ws.Range("I107").Value = 646.1667
$ws.Range("K107").Value = 646.1667
$ws.Range("M107").Value = 1273.8333

$ws.Range("H125").Value = 5185.909
$ws.Range("I125").Value = 288
$ws.Range("J125").Value = 6274.3335
$ws.Range("K125").Value = 2592
$ws.Range("L125").Value = 56469.0015
$ws.Range("M125").Value = -132
$ws.Range("N125").Value = -61389.0015

$ws.Range("H132").Value = 1249.2458
$ws.Range("I132").Value = 1284.375
$ws.Range("J132").Value = 855.8
$ws.Range("K132").Value = 3853.125
$ws.Range("L132").Value = 2567.4
$ws.Range("M132").Value = -1323.125
$ws.Range("N132").Value = -7627.4

$ws.Range("H137").Value = 2118.4746
$ws.Range("I137").Value = 1912.6666
$ws.Range("J137").Value = 2626.9412
$ws.Range("K137").Value = 5737.9998
$ws.Range("L137").Value = 7880.823600000001
$ws.Range("M137").Value = -3187.9998
$ws.Range("N137").Value = -12980.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19425.074
$ws.Range("I32").Value = 20291.39
$ws.Range("K32").Value = 20291.39
$ws.Range("M32").Value = -20004.39

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 67795
$ws.Range("J13").Value = 67795
$ws.Range("L13").Value = 67795
$ws.Range("N13").Value = -68131

$ws.Range("H99").Value = 830.3077
$ws.Range("I99").Value = 1250
$ws.Range("J99").Value = 643.7778
$ws.Range("K99").Value = 1250
$ws.Range("L99").Value = 643.7778
$ws.Range("M99").Value = 248
$ws.Range("N99").Value = -3639.7778

$ws.Range("H134").Value = 2215.641
$ws.Range("I134").Value = 2156.7144
$ws.Range("J134").Value = 2365.6365
$ws.Range("K134").Value = 6470.1432
$ws.Range("L134").Value = 7096.9095
$ws.Range("M134").Value = -3935.1432
$ws.Range("N134").Value = -12166.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6035.613
$ws.Range("I31").Value = 9019.929
$ws.Range("J31").Value = 3577.9412
$ws.Range("K31").Value = 9019.929
$ws.Range("L31").Value = 3577.9412
$ws.Range("M31").Value = -8724.929
$ws.Range("N31").Value = -4167.9412

$ws.Range("H34").Value = 6035.613
$ws.Range("I34").Value = 9019.929
$ws.Range("J34").Value = 3577.9412
$ws.Range("K34").Value = 9019.929
$ws.Range("L34").Value = 3577.9412
$ws.Range("M34").Value = -8817.929
$ws.Range("N34").Value = -3981.9412

$ws.Range("H53").Value = 34460
$ws.Range("J53").Value = 34460
$ws.Range("L53").Value = 34460
$ws.Range("N53").Value = -35674

$ws.Range("H58").Value = 1492016.6
$ws.Range("I58").Value = 2218425.2
$ws.Range("J58").Value = 2878.75
$ws.Range("K58").Value = 2218425.2
$ws.Range("L58").Value = 2878.75
$ws.Range("M58").Value = -2218222.2
$ws.Range("N58").Value = -3284.75

$ws.Range("H132").Value = 2097.3062
$ws.Range("I132").Value = 1758.8125
$ws.Range("J132").Value = 2734.4707
$ws.Range("K132").Value = 5276.4375
$ws.Range("L132").Value = 8203.4121
$ws.Range("M132").Value = -2746.4375
$ws.Range("N132").Value = -13263.4121

$ws.Range("H136").Value = 1492016.6
$ws.Range("I136").Value = 2218425.2
$ws.Range("J136").Value = 2878.75
$ws.Range("K136").Value = 6655275.600000001
$ws.Range("L136").Value = 8636.25
$ws.Range("M136").Value = -6652725.600000001
$ws.Range("N136").Value = -13736.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3223.9375
$ws.Range("J127").Value = 3223.9375
$ws.Range("L127").Value = 9671.8125
$ws.Range("N127").Value = -19591.8125

$ws.Range("H131").Value = 65570.60000000001
$ws.Range("I131").Value = 2971.6667
$ws.Range("J131").Value = 107303.22
$ws.Range("K131").Value = 8915.000100000001
$ws.Range("L131").Value = 321909.66
$ws.Range("M131").Value = -3875.000100000001
$ws.Range("N131").Value = -331989.66

$ws.Range("H132").Value = 1204.175
$ws.Range("I132").Value = 1703.0834
$ws.Range("J132").Value = 990.3570999999999
$ws.Range("K132").Value = 15327.7506
$ws.Range("L132").Value = 8913.213899999999
$ws.Range("M132").Value = -12797.7506
$ws.Range("N132").Value = -13973.2139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1074.7826
$ws.Range("I97").Value = 1011.3333
$ws.Range("K97").Value = 1011.3333
$ws.Range("M97").Value = -515.3333

$ws.Range("H122").Value = 11070.615
$ws.Range("I122").Value = 13599.875
$ws.Range("J122").Value = 7023.8
$ws.Range("K122").Value = 40799.625
$ws.Range("L122").Value = 21071.4
$ws.Range("M122").Value = -38349.625
$ws.Range("N122").Value = -25971.4

$ws.Range("H126").Value = 2986.25
$ws.Range("I126").Value = 1925
$ws.Range("J126").Value = 3693.75
$ws.Range("K126").Value = 5775
$ws.Range("L126").Value = 11081.25
$ws.Range("M126").Value = -3305
$ws.Range("N126").Value = -16021.25

$ws.Range("H132").Value = 1771.2778
$ws.Range("I132").Value = 1625.6666
$ws.Range("J132").Value = 2936.1667
$ws.Range("K132").Value = 4876.9998
$ws.Range("L132").Value = 8808.500100000001
$ws.Range("M132").Value = -2346.9998
$ws.Range("N132").Value = -13868.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 250.65
$ws.Range("I93").Value = 209.83333
$ws.Range("K93").Value = 209.83333
$ws.Range("M93").Value = 1038.16667

$ws.Range("H117").Value = 49028.57
$ws.Range("J117").Value = 49028.57
$ws.Range("L117").Value = 49028.57
$ws.Range("N117").Value = -58206.57

$ws.Range("H136").Value = 4125.415
$ws.Range("I136").Value = 2107.9143
$ws.Range("J136").Value = 8048.3335
$ws.Range("K136").Value = 6323.742899999999
$ws.Range("L136").Value = 24145.0005
$ws.Range("M136").Value = -3773.742899999999
$ws.Range("N136").Value = -29245.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2434.2222
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 2638.5
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 7915.5
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -11755.5

$ws.Range("H118").Value = 67500
$ws.Range("J118").Value = 67500
$ws.Range("L118").Value = 67500
$ws.Range("N118").Value = -70814

$ws.Range("H132").Value = 1600.8959
$ws.Range("I132").Value = 895.8125
$ws.Range("K132").Value = 2687.4375
$ws.Range("M132").Value = -157.4375

$ws.Range("H136").Value = 6532.84
$ws.Range("I136").Value = 4750.7188
$ws.Range("J136").Value = 9701.056
$ws.Range("K136").Value = 14252.1564
$ws.Range("L136").Value = 29103.168
$ws.Range("M136").Value = -11702.1564
$ws.Range("N136").Value = -34203.16800000001

